$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44243
$ws.Range("K2").Value = 11000
$ws.Range("L2").Value = 12000
$ws.Range("M2").Value = 11500
$ws.Range("P2").Value = 192

# Row 3
$ws.Range("D3").Value = 44203
$ws.Range("K3").Value = 11000
$ws.Range("L3").Value = 12000
$ws.Range("M3").Value = 11500
$ws.Range("N3").Value = '$/caja 60 unidades'
$ws.Range("P3").Value = 192
$ws.Range("Q3").Value = 60

# Row 4
$ws.Range("D4").Value = 44565
$ws.Range("J4").Value = 100
$ws.Range("K4").Value = 7000
$ws.Range("L4").Value = 8000
$ws.Range("M4").Value = 7500
$ws.Range("O4").Value = 'Región de Arica y Parinacota'
$ws.Range("P4").Value = 125

# Row 5
$ws.Range("D5").Value = 44377
$ws.Range("K5").Value = 11000
$ws.Range("L5").Value = 12000
$ws.Range("M5").Value = 11500
$ws.Range("P5").Value = 192

# Row 6
$ws.Range("D6").Value = 44568
$ws.Range("K6").Value = 7000
$ws.Range("L6").Value = 7500
$ws.Range("M6").Value = 7250
$ws.Range("P6").Value = 121

# Row 7
$ws.Range("D7").Value = 44355
$ws.Range("K7").Value = 10000
$ws.Range("L7").Value = 11000
$ws.Range("M7").Value = 10500
$ws.Range("P7").Value = 175

# Row 8
$ws.Range("D8").Value = 44432
$ws.Range("K8").Value = 14000
$ws.Range("L8").Value = 15000
$ws.Range("M8").Value = 14500
$ws.Range("P8").Value = 242

# Row 9
$ws.Range("D9").Value = 44323
$ws.Range("K9").Value = 9000
$ws.Range("L9").Value = 10000
$ws.Range("M9").Value = 9500
$ws.Range("P9").Value = 158

# Row 10
$ws.Range("D10").Value = 44474
$ws.Range("J10").Value = 100
$ws.Range("K10").Value = 19000
$ws.Range("L10").Value = 20000
$ws.Range("M10").Value = 19500
$ws.Range("P10").Value = 325

# Row 11
$ws.Range("D11").Value = 44238
$ws.Range("K11").Value = 12000
$ws.Range("L11").Value = 14000
$ws.Range("M11").Value = 13000
$ws.Range("P11").Value = 217

# Row 12
$ws.Range("D12").Value = 44385
$ws.Range("K12").Value = 15000
$ws.Range("L12").Value = 16000
$ws.Range("M12").Value = 15500
$ws.Range("P12").Value = 258

# Row 13
$ws.Range("D13").Value = 44447
$ws.Range("J13").Value = 100
$ws.Range("K13").Value = 16000
$ws.Range("L13").Value = 17000
$ws.Range("M13").Value = 16500
$ws.Range("P13").Value = 275

# Row 14
$ws.Range("D14").Value = 44469

# Row 15
$ws.Range("D15").Value = 44281
$ws.Range("K15").Value = 12000
$ws.Range("L15").Value = 13000
$ws.Range("M15").Value = 12500
$ws.Range("P15").Value = 208

# Row 17
$ws.Range("D17").Value = 44370
$ws.Range("K17").Value = 15000
$ws.Range("L17").Value = 16000
$ws.Range("M17").Value = 15500
$ws.Range("O17").Value = 'Región Metropolitana'
$ws.Range("P17").Value = 258

# Row 18
$ws.Range("D18").Value = 44551
$ws.Range("J18").Value = 100
$ws.Range("K18").Value = 7000
$ws.Range("L18").Value = 8000
$ws.Range("M18").Value = 7500
$ws.Range("O18").Value = 'Región de Arica y Parinacota'
$ws.Range("P18").Value = 125

# Row 19
$ws.Range("D19").Value = 44484
$ws.Range("J19").Value = 450
$ws.Range("K19").Value = 11000
$ws.Range("L19").Value = 12000
$ws.Range("M19").Value = 11556
$ws.Range("P19").Value = 193

# Row 20
$ws.Range("D20").Value = 44383
$ws.Range("J20").Value = 100
$ws.Range("K20").Value = 14000
$ws.Range("L20").Value = 15000
$ws.Range("M20").Value = 14500
$ws.Range("O20").Value = 'Región de Arica y Parinacota'
$ws.Range("P20").Value = 242

# Row 21
$ws.Range("D21").Value = 44358
$ws.Range("J21").Value = 100
$ws.Range("K21").Value = 11000
$ws.Range("L21").Value = 12000
$ws.Range("M21").Value = 11500
$ws.Range("P21").Value = 192

# Row 22
$ws.Range("D22").Value = 44425
$ws.Range("K22").Value = 16000
$ws.Range("L22").Value = 17000
$ws.Range("M22").Value = 16500
$ws.Range("P22").Value = 275

# Row 23
$ws.Range("D23").Value = 44365
$ws.Range("K23").Value = 13000
$ws.Range("L23").Value = 14000
$ws.Range("M23").Value = 13500
$ws.Range("P23").Value = 225

# Row 24
$ws.Range("D24").Value = 44253
$ws.Range("J24").Value = 200
$ws.Range("K24").Value = 9000
$ws.Range("L24").Value = 10000
$ws.Range("M24").Value = 9500
$ws.Range("P24").Value = 158

# Row 25
$ws.Range("D25").Value = 44217
$ws.Range("J25").Value = 200
$ws.Range("K25").Value = 8000
$ws.Range("L25").Value = 9000
$ws.Range("M25").Value = 8500
$ws.Range("O25").Value = 'Región del Maule'
$ws.Range("P25").Value = 142

# Row 26
$ws.Range("D26").Value = 44537
$ws.Range("J26").Value = 220
$ws.Range("K26").Value = 9000
$ws.Range("L26").Value = 10000
$ws.Range("M26").Value = 9545
$ws.Range("P26").Value = 159

# Row 27
$ws.Range("D27").Value = 44420
$ws.Range("J27").Value = 200
$ws.Range("K27").Value = 16000
$ws.Range("L27").Value = 17000
$ws.Range("M27").Value = 16500
$ws.Range("P27").Value = 275

# Row 28
$ws.Range("D28").Value = 44246
$ws.Range("J28").Value = 200
$ws.Range("K28").Value = 10000
$ws.Range("L28").Value = 12000
$ws.Range("M28").Value = 11000
$ws.Range("O28").Value = 'Región del Maule'
$ws.Range("P28").Value = 183

# Row 29
$ws.Range("D29").Value = 44399
$ws.Range("K29").Value = 16000
$ws.Range("L29").Value = 17000
$ws.Range("M29").Value = 16500
$ws.Range("P29").Value = 275

# Row 30
$ws.Range("D30").Value = 44258
$ws.Range("J30").Value = 200
$ws.Range("K30").Value = 12000
$ws.Range("L30").Value = 13000
$ws.Range("M30").Value = 12500
$ws.Range("P30").Value = 208

# Row 31
$ws.Range("D31").Value = 44572
$ws.Range("J31").Value = 310
$ws.Range("K31").Value = 5500
$ws.Range("L31").Value = 6000
$ws.Range("M31").Value = 5742
$ws.Range("P31").Value = 96

# Row 32
$ws.Range("D32").Value = 44519
$ws.Range("J32").Value = 450
$ws.Range("K32").Value = 6500
$ws.Range("L32").Value = 7000
$ws.Range("M32").Value = 6778
$ws.Range("P32").Value = 113

# Row 33
$ws.Range("D33").Value = 44223
$ws.Range("J33").Value = 200
$ws.Range("K33").Value = 9000
$ws.Range("L33").Value = 10000
$ws.Range("M33").Value = 9500
$ws.Range("P33").Value = 158

# Row 34
$ws.Range("D34").Value = 44189
$ws.Range("J34").Value = 100
$ws.Range("M34").Value = 11500
$ws.Range("P34").Value = 192

# Row 35
$ws.Range("D35").Value = 44285
$ws.Range("K35").Value = 12000
$ws.Range("L35").Value = 13000
$ws.Range("M35").Value = 12500
$ws.Range("P35").Value = 208

# Row 36
$ws.Range("D36").Value = 44455
$ws.Range("J36").Value = 100
$ws.Range("K36").Value = 15000
$ws.Range("L36").Value = 16000
$ws.Range("M36").Value = 15500
$ws.Range("N36").Value = '$/caja 50 unidades'
$ws.Range("O36").Value = 'Región de Arica y Parinacota'
$ws.Range("P36").Value = 310
$ws.Range("Q36").Value = 50

# Row 37
$ws.Range("D37").Value = 44397
$ws.Range("K37").Value = 17000
$ws.Range("L37").Value = 18000
$ws.Range("M37").Value = 17500
$ws.Range("P37").Value = 292

# Row 38
$ws.Range("D38").Value = 44336
$ws.Range("J38").Value = 100
$ws.Range("K38").Value = 10000
$ws.Range("L38").Value = 11000
$ws.Range("M38").Value = 10500
$ws.Range("P38").Value = 175

# Row 39
$ws.Range("D39").Value = 44166
$ws.Range("K39").Value = 6500
$ws.Range("L39").Value = 7000
$ws.Range("M39").Value = 6750
$ws.Range("P39").Value = 112

# Row 40
$ws.Range("D40").Value = 44435

# Row 41
$ws.Range("D41").Value = 44530
$ws.Range("J41").Value = 350
$ws.Range("K41").Value = 6000
$ws.Range("L41").Value = 6500
$ws.Range("M41").Value = 6286
$ws.Range("P41").Value = 79

# Row 42
$ws.Range("D42").Value = 44187
$ws.Range("J42").Value = 200
$ws.Range("K42").Value = 8000
$ws.Range("L42").Value = 9000
$ws.Range("M42").Value = 8500
$ws.Range("P42").Value = 142

# Row 43
$ws.Range("D43").Value = 44567
$ws.Range("K43").Value = 7000
$ws.Range("L43").Value = 7500
$ws.Range("M43").Value = 7250
$ws.Range("P43").Value = 121

# Row 44
$ws.Range("D44").Value = 44237
$ws.Range("K44").Value = 10000
$ws.Range("L44").Value = 11000
$ws.Range("M44").Value = 10500
$ws.Range("P44").Value = 175

# Row 45
$ws.Range("D45").Value = 44421
$ws.Range("J45").Value = 100
$ws.Range("K45").Value = 17000
$ws.Range("L45").Value = 18000
$ws.Range("M45").Value = 17500
$ws.Range("O45").Value = 'Región de Arica y Parinacota'
$ws.Range("P45").Value = 292

# Row 46
$ws.Range("D46").Value = 44350
$ws.Range("K46").Value = 10000
$ws.Range("M46").Value = 11000
$ws.Range("P46").Value = 183

# Row 47
$ws.Range("D47").Value = 44453
$ws.Range("K47").Value = 16000
$ws.Range("L47").Value = 17000
$ws.Range("M47").Value = 16500
$ws.Range("P47").Value = 275

# Row 48
$ws.Range("D48").Value = 44215
$ws.Range("K48").Value = 9000
$ws.Range("L48").Value = 10000
$ws.Range("M48").Value = 9500
$ws.Range("P48").Value = 158

# Row 49
$ws.Range("D49").Value = 44539
$ws.Range("J49").Value = 180
$ws.Range("M49").Value = 6722

# Row 50
$ws.Range("D50").Value = 44523
$ws.Range("J50").Value = 100
$ws.Range("K50").Value = 6500
$ws.Range("L50").Value = 7000
$ws.Range("M50").Value = 6750
$ws.Range("O50").Value = 'Región de Arica y Parinacota'
$ws.Range("P50").Value = 112

# Row 51
$ws.Range("D51").Value = 44372
$ws.Range("K51").Value = 14000
$ws.Range("L51").Value = 15000
$ws.Range("M51").Value = 14500
$ws.Range("P51").Value = 242

# Row 52
$ws.Range("D52").Value = 44292
$ws.Range("K52").Value = 14000
$ws.Range("L52").Value = 15000
$ws.Range("M52").Value = 14500
$ws.Range("P52").Value = 242

# Row 53
$ws.Range("D53").Value = 44505
$ws.Range("J53").Value = 300
$ws.Range("K53").Value = 6500
$ws.Range("L53").Value = 7000
$ws.Range("M53").Value = 6750
$ws.Range("N53").Value = '$/caja 80 unidades'
$ws.Range("O53").Value = 'Región del Maule'
$ws.Range("P53").Value = 84
$ws.Range("Q53").Value = 80

# Row 54
$ws.Range("D54").Value = 44168
$ws.Range("J54").Value = 100
$ws.Range("K54").Value = 6500
$ws.Range("L54").Value = 7000
$ws.Range("M54").Value = 6750
$ws.Range("P54").Value = 112

# Row 55
$ws.Range("D55").Value = 44299
$ws.Range("K55").Value = 14000
$ws.Range("L55").Value = 15000
$ws.Range("M55").Value = 14500
$ws.Range("P55").Value = 242

# Row 56
$ws.Range("D56").Value = 44498
$ws.Range("J56").Value = 350
$ws.Range("K56").Value = 6500
$ws.Range("L56").Value = 7000
$ws.Range("M56").Value = 6786
$ws.Range("P56").Value = 113

# Row 57
$ws.Range("D57").Value = 44308
$ws.Range("K57").Value = 11000
$ws.Range("L57").Value = 12000
$ws.Range("M57").Value = 11500
$ws.Range("P57").Value = 192

# Row 58
$ws.Range("D58").Value = 44357

# Row 59
$ws.Range("D59").Value = 44320
$ws.Range("J59").Value = 100

# Row 60
$ws.Range("D60").Value = 44306
$ws.Range("J60").Value = 200
$ws.Range("K60").Value = 9000
$ws.Range("L60").Value = 10000
$ws.Range("M60").Value = 9500
$ws.Range("P60").Value = 158

# Row 61
$ws.Range("D61").Value = 44295
$ws.Range("K61").Value = 13000
$ws.Range("L61").Value = 14000
$ws.Range("M61").Value = 13500
$ws.Range("P61").Value = 225

# Row 62
$ws.Range("D62").Value = 44210
$ws.Range("J62").Value = 100
$ws.Range("K62").Value = 11000
$ws.Range("L62").Value = 12000
$ws.Range("M62").Value = 11500
$ws.Range("P62").Value = 192

# Row 63
$ws.Range("D63").Value = 44343
$ws.Range("K63").Value = 10000
$ws.Range("L63").Value = 11000
$ws.Range("M63").Value = 10500
$ws.Range("O63").Value = 'Región de Arica y Parinacota'
$ws.Range("P63").Value = 175

# Row 64
$ws.Range("D64").Value = 44230
$ws.Range("K64").Value = 9000
$ws.Range("L64").Value = 10000
$ws.Range("M64").Value = 9500
$ws.Range("P64").Value = 158

# Row 65
$ws.Range("D65").Value = 44316
$ws.Range("J65").Value = 100
$ws.Range("K65").Value = 9000
$ws.Range("L65").Value = 10000
$ws.Range("M65").Value = 9500
$ws.Range("P65").Value = 158

# Row 66
$ws.Range("D66").Value = 44265
$ws.Range("J66").Value = 100
$ws.Range("K66").Value = 13000
$ws.Range("L66").Value = 15000
$ws.Range("M66").Value = 14000
$ws.Range("P66").Value = 233

# Row 67
$ws.Range("D67").Value = 44460
$ws.Range("J67").Value = 100
$ws.Range("K67").Value = 16000
$ws.Range("L67").Value = 17000
$ws.Range("M67").Value = 16500
$ws.Range("P67").Value = 275

# Row 68
$ws.Range("D68").Value = 44526
$ws.Range("J68").Value = 200
$ws.Range("K68").Value = 7000
$ws.Range("L68").Value = 8000
$ws.Range("M68").Value = 7500
$ws.Range("P68").Value = 125

# Row 69
$ws.Range("D69").Value = 44272
$ws.Range("J69").Value = 100
$ws.Range("K69").Value = 12000
$ws.Range("L69").Value = 13000
$ws.Range("M69").Value = 12500
$ws.Range("P69").Value = 208

# Row 70
$ws.Range("D70").Value = 44467
$ws.Range("K70").Value = 15000
$ws.Range("L70").Value = 16000
$ws.Range("M70").Value = 15500
$ws.Range("P70").Value = 258

# Row 71
$ws.Range("D71").Value = 44211
$ws.Range("J71").Value = 200
$ws.Range("K71").Value = 11000
$ws.Range("L71").Value = 12000
$ws.Range("M71").Value = 11500
$ws.Range("P71").Value = 192

# Row 72
$ws.Range("D72").Value = 44313
$ws.Range("K72").Value = 9000
$ws.Range("L72").Value = 10000
$ws.Range("M72").Value = 9500
$ws.Range("P72").Value = 158

# Row 73
$ws.Range("D73").Value = 44334
$ws.Range("K73").Value = 11000
$ws.Range("L73").Value = 12000
$ws.Range("M73").Value = 11500
$ws.Range("P73").Value = 192

# Row 74
$ws.Range("D74").Value = 44517
$ws.Range("J74").Value = 250
$ws.Range("K74").Value = 5500
$ws.Range("L74").Value = 6000
$ws.Range("M74").Value = 5700
$ws.Range("P74").Value = 95

# Row 75
$ws.Range("D75").Value = 44330
$ws.Range("K75").Value = 10000
$ws.Range("L75").Value = 11000
$ws.Range("M75").Value = 10500
$ws.Range("P75").Value = 175

# Row 76
$ws.Range("D76").Value = 44196
$ws.Range("K76").Value = 11000
$ws.Range("L76").Value = 12000
$ws.Range("M76").Value = 11500
$ws.Range("P76").Value = 192

# Row 77
$ws.Range("D77").Value = 44463
$ws.Range("K77").Value = 15000
$ws.Range("L77").Value = 15500
$ws.Range("M77").Value = 15250
$ws.Range("P77").Value = 254

# Row 78
$ws.Range("D78").Value = 44476
$ws.Range("K78").Value = 16000
$ws.Range("L78").Value = 17000
$ws.Range("M78").Value = 16500
$ws.Range("P78").Value = 275

# Row 79
$ws.Range("D79").Value = 44496
$ws.Range("J79").Value = 350
$ws.Range("K79").Value = 6500
$ws.Range("L79").Value = 7000
$ws.Range("M79").Value = 6786
$ws.Range("O79").Value = 'Región de Coquimbo'
$ws.Range("P79").Value = 113

# Row 80
$ws.Range("D80").Value = 44250
$ws.Range("J80").Value = 100
$ws.Range("K80").Value = 8000
$ws.Range("L80").Value = 9000
$ws.Range("M80").Value = 8500
$ws.Range("P80").Value = 142

# Row 81
$ws.Range("D81").Value = 44552
$ws.Range("J81").Value = 100
$ws.Range("K81").Value = 7000
$ws.Range("L81").Value = 8000
$ws.Range("M81").Value = 7500
$ws.Range("P81").Value = 125

# Row 82
$ws.Range("D82").Value = 44379
$ws.Range("K82").Value = 13000
$ws.Range("L82").Value = 14000
$ws.Range("M82").Value = 13500
$ws.Range("P82").Value = 225

# Row 83
$ws.Range("D83").Value = 44509
$ws.Range("K83").Value = 6000
$ws.Range("L83").Value = 6500
$ws.Range("M83").Value = 6250
$ws.Range("P83").Value = 104

# Row 84
$ws.Range("D84").Value = 44488
$ws.Range("K84").Value = 7000
$ws.Range("L84").Value = 7500
$ws.Range("M84").Value = 7250
$ws.Range("P84").Value = 121

# Row 85
$ws.Range("D85").Value = 44341
$ws.Range("K85").Value = 9000
$ws.Range("L85").Value = 10000
$ws.Range("M85").Value = 9500
$ws.Range("P85").Value = 158

# Row 86
$ws.Range("D86").Value = 44482
$ws.Range("J86").Value = 350
$ws.Range("K86").Value = 10000
$ws.Range("L86").Value = 11000
$ws.Range("M86").Value = 10429
$ws.Range("P86").Value = 174

# Row 87
$ws.Range("D87").Value = 44294
$ws.Range("K87").Value = 14000
$ws.Range("L87").Value = 15000
$ws.Range("M87").Value = 14500
$ws.Range("P87").Value = 242

# Row 88
$ws.Range("D88").Value = 44558
$ws.Range("J88").Value = 250
$ws.Range("K88").Value = 8000
$ws.Range("L88").Value = 8500
$ws.Range("M88").Value = 8300
$ws.Range("O88").Value = 'Provincia de Limarí'
$ws.Range("P88").Value = 138

# Row 89
$ws.Range("D89").Value = 44278
$ws.Range("J89").Value = 100
$ws.Range("K89").Value = 11000
$ws.Range("L89").Value = 12000
$ws.Range("M89").Value = 11500
$ws.Range("P89").Value = 192

# Row 90
$ws.Range("D90").Value = 44392
$ws.Range("K90").Value = 16000
$ws.Range("L90").Value = 17000
$ws.Range("M90").Value = 16500
$ws.Range("P90").Value = 275

# Row 91
$ws.Range("D91").Value = 44194
$ws.Range("J91").Value = 100
$ws.Range("K91").Value = 11000
$ws.Range("L91").Value = 12000
$ws.Range("M91").Value = 11500
$ws.Range("P91").Value = 192

# Row 92
$ws.Range("D92").Value = 44162
$ws.Range("J92").Value = 200
$ws.Range("K92").Value = 7000
$ws.Range("L92").Value = 7500
$ws.Range("M92").Value = 7250
$ws.Range("O92").Value = 'Región del Maule'
$ws.Range("P92").Value = 121

# Row 93
$ws.Range("D93").Value = 44532
$ws.Range("J93").Value = 250
$ws.Range("K93").Value = 6500
$ws.Range("L93").Value = 7000
$ws.Range("M93").Value = 6700
$ws.Range("P93").Value = 112

# Row 94
$ws.Range("D94").Value = 44427
$ws.Range("K94").Value = 14000
$ws.Range("L94").Value = 15000
$ws.Range("M94").Value = 14500
$ws.Range("P94").Value = 242

# Row 95
$ws.Range("D95").Value = 44491
$ws.Range("K95").Value = 8500
$ws.Range("L95").Value = 9000
$ws.Range("M95").Value = 8750
$ws.Range("P95").Value = 146

# Row 96
$ws.Range("D96").Value = 44390
$ws.Range("K96").Value = 16000
$ws.Range("L96").Value = 17000
$ws.Range("M96").Value = 16500
$ws.Range("P96").Value = 275

# Row 97
$ws.Range("D97").Value = 44327
$ws.Range("K97").Value = 9000
$ws.Range("L97").Value = 10000
$ws.Range("M97").Value = 9500
$ws.Range("P97").Value = 158

# Row 98
$ws.Range("D98").Value = 44503
$ws.Range("J98").Value = 250
$ws.Range("K98").Value = 7500
$ws.Range("L98").Value = 8000
$ws.Range("M98").Value = 7700
$ws.Range("P98").Value = 128

# Row 99
$ws.Range("D99").Value = 44462
$ws.Range("K99").Value = 14500
$ws.Range("L99").Value = 15000
$ws.Range("M99").Value = 14750
$ws.Range("P99").Value = 246

# Row 100
$ws.Range("D100").Value = 44511
$ws.Range("K100").Value = 8000
$ws.Range("L100").Value = 9000
$ws.Range("M100").Value = 8500
$ws.Range("P100").Value = 142
